# Applies the cell-level text/value updates from the commit diff
# ("Updated cryptos list ... with GitHub Actions") to cryptos.xlsx.
# Column D holds numeric-looking price strings that must remain exact
# text (decimal places / thousands separators matter), so for any new
# value that Excel would otherwise auto-convert to a number we force
# the cell to Text format (NumberFormat '@') immediately before writing it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.388.93'
$ws.Range('E2').Value = '  -0.71%  '
$ws.Range('D3').Value = '1.848.43'
$ws.Range('E3').Value = '  -0.33%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9990'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6345'
$ws.Range('E6').Value = '  -0.99%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07572'
$ws.Range('E8').Value = '  -0.26%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2963'
$ws.Range('E9').Value = '  -1.39%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.79'
$ws.Range('E10').Value = '  +1.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07731'
$ws.Range('E11').Value = '  +0.79%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.012'
$ws.Range('E12').Value = '  -0.91%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6832'
$ws.Range('B14').Value = 'Litecoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '83.06'
$ws.Range('E14').Value = '  -1.27%  '
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000009975'
$ws.Range('E15').Value = '  +3.04%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.136'
$ws.Range('E16').Value = '  -2.50%  '
$ws.Range('D17').Value = '29.415.19'
$ws.Range('E17').Value = '  -0.80%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '231.23'
$ws.Range('E18').Value = '  -3.29%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.556'
$ws.Range('E21').Value = '  -1.26%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('E23').Value = '  +228.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '16.66'
$ws.Range('E24').Value = '  +171.56%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '156.42'
$ws.Range('E25').Value = '  -0.64%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1398'
$ws.Range('E26').Value = '  -0.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.386'
$ws.Range('E27').Value = '  -1.59%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.720'
$ws.Range('E29').Value = '  +172.39%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.472'
$ws.Range('E30').Value = '  -1.32%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.05695'
$ws.Range('E31').Value = '  -3.57%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.255'
$ws.Range('E32').Value = '  -2.18%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.133'
$ws.Range('E33').Value = '  -0.47%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.027'
$ws.Range('E34').Value = '  -1.67%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.850'
$ws.Range('E35').Value = '  -3.40%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.157'
$ws.Range('E36').Value = '  -2.53%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.7167'
$ws.Range('E37').Value = '  -1.35%  '
$ws.Range('E38').Value = '  -0.08%  '
$ws.Range('D39').Value = '1.244.00'
$ws.Range('E39').Value = '  +2.22%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.803'
$ws.Range('E40').Value = '  -0.16%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.01811'
$ws.Range('E41').Value = '  +1.72%  '
$ws.Range('E42').Value = '  +264.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9041'
$ws.Range('E43').Value = '  -1.47%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.0000'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '101.88'
$ws.Range('E45').Value = '  -0.16%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '66.22'
$ws.Range('E46').Value = '  -1.87%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.076'
$ws.Range('E47').Value = '  -5.70%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.126'
$ws.Range('E48').Value = '  -0.88%  '
$ws.Range('B49').Value = 'TheSandbox'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4025'
$ws.Range('E49').Value = '  -1.19%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.710'
$ws.Range('E50').Value = '  +1.28%  '
$ws.Range('E51').Value = '  -0.31%  '
